$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Silan%Gu%NULL%1,         Yanfei%Chen%NULL%1,         Zhengjie%Wu%NULL%1,         Yunbo%Chen%NULL%1,         Hainv%Gao%NULL%1,         Longxian%Lv%NULL%1,         Feifei%Guo%NULL%1,         Xuewu%Zhang%NULL%1,         Rui%Luo%NULL%1,         Chenjie%Huang%NULL%1,         Haifeng%Lu%NULL%1,         Beiwen%Zheng%NULL%1,         Jiaying%Zhang%NULL%1,         Ren%Yan%NULL%1,         Hua%Zhang%NULL%0,         Huiyong%Jiang%NULL%1,         Qiaomai%Xu%NULL%1,         Jing%Guo%NULL%1,         Yiwen%Gong%NULL%1,         Lingling%Tang%NULL%1,         Lanjuan%Li%ljli@zju.edu.cn%0]'
$ws.Range("E3").Value = '[Tao%Zuo%NULL%2,         Fen%Zhang%NULL%2,         Grace C.Y.%Lui%NULL%2,         Yun Kit%Yeoh%NULL%1,         Amy Y.L.%Li%NULL%1,         Hui%Zhan%NULL%2,         Yating%Wan%NULL%1,         Arthur C.K.%Chung%NULL%1,         Chun Pan%Cheung%NULL%1,         Nan%Chen%NULL%2,         Christopher K.C.%Lai%NULL%1,         Zigui%Chen%NULL%1,         Eugene Y.K.%Tso%NULL%2,         Kitty S.C.%Fung%NULL%1,         Veronica%Chan%NULL%1,         Lowell%Ling%NULL%1,         Gavin%Joynt%NULL%1,         David S.C.%Hui%NULL%0,         Francis K.L.%Chan%NULL%2,         Paul K.S.%Chan%NULL%2,         Siew C.%Ng%NULL%2]'
$ws.Range("E4").Value = '[Tao%Zuo%NULL%0,         Hui%Zhan%NULL%0,         Fen%Zhang%NULL%0,         Qin%Liu%NULL%0,         Eugene Y.K.%Tso%NULL%0,         Grace C.Y.%Lui%NULL%0,         Nan%Chen%NULL%0,         Amy%Li%NULL%1,         Wenqi%Lu%NULL%1,         Francis K.L.%Chan%NULL%0,         Paul K.S.%Chan%NULL%0,         Siew C.%Ng%NULL%0]'
$ws.Range("E5").Value = '[Haocheng%Zhang%NULL%1,         Jing-Wen%Ai%NULL%1,         Wenjiao%Yang%NULL%1,         Xian%Zhou%NULL%1,         Fusheng%He%NULL%1,         Shumei%Xie%NULL%1,         Weiqi%Zeng%NULL%1,         Yang%Li%NULL%0,         Yiqi%Yu%NULL%1,         Xuejing%Gou%NULL%1,         Yongjun%Li%NULL%1,         Xiaorui%Wang%NULL%1,         Hang%Su%NULL%1,         Teng%Xu%NULL%2,         Wenhong%Zhang%wenhongzhang_hs@126.com%0]'
$ws.Range("E6").Value = '[Flavio%De Maio%NULL%1,         Brunella%Posteraro%NULL%1,         Francesca Romana%Ponziani%NULL%1,         Paola%Cattani%NULL%1,         Antonio%Gasbarrini%NULL%0,         Maurizio%Sanguinetti%maurizio.sanguinetti@unicatt.it%1]'
$ws.Range("E7").Value = '[Liangjun%Chen%NULL%1,         Weiyong%Liu%NULL%0,         Qi%Zhang%NULL%1,         Ke%Xu%NULL%2,         Ke%Xu%NULL%0,         Guangming%Ye%NULL%0,         Weichen%Wu%NULL%1,         Ziyong%Sun%NULL%1,         Fang%Liu%NULL%0,         Kailang%Wu%NULL%0,         Bo%Zhong%NULL%1,         Yi%Mei%NULL%1,         Wenxia%Zhang%NULL%1,         Yu%Chen%NULL%0,         Yirong%Li%NULL%0,         Yirong%Li%NULL%0,         Mang%Shi%NULL%1,         Ke%Lan%NULL%1,         Yingle%Liu%NULL%1]'
$ws.Range("E8").Value = '[Li-Li%Ren%NULL%0,         Ye-Ming%Wang%NULL%1,         Zhi-Qiang%Wu%NULL%1,         Zi-Chun%Xiang%NULL%1,         Li%Guo%NULL%0,         Teng%Xu%NULL%0,         Yong-Zhong%Jiang%NULL%1,         Yan%Xiong%NULL%0,         Yong-Jun%Li%NULL%1,         Xing-Wang%Li%NULL%1,         Hui%Li%NULL%0,         Guo-Hui%Fan%NULL%1,         Xiao-Ying%Gu%NULL%1,         Yan%Xiao%NULL%0,         Hong%Gao%NULL%0,         Jiu-Yang%Xu%NULL%1,         Fan%Yang%NULL%0,         Xin-Ming%Wang%NULL%1,         Chao%Wu%NULL%1,         Lan%Chen%NULL%1,         Yi-Wei%Liu%NULL%1,         Bo%Liu%NULL%0,         Jian%Yang%NULL%1,         Xiao-Rui%Wang%NULL%1,         Jie%Dong%NULL%1,         Li%Li%NULL%0,         Chao-Lin%Huang%NULL%1,         Jian-Ping%Zhao%NULL%1,         Yi%Hu%NULL%0,         Zhen-Shun%Cheng%NULL%1,         Lin-Lin%Liu%NULL%1,         Zhao-Hui%Qian%NULL%1,         Chuan%Qin%NULL%0,         Qi%Jin%NULL%0,         Bin%Cao%NULL%0,         Jian-Wei%Wang%NULL%1,         Xiu-Yuan%Hao%NULL%0,         Xiu-Yuan%Hao%NULL%0,         Pei-Fang%Wei%NULL%0]'
$ws.Range("E9").Value = '[Jun%Fan%NULL%1,         Xiang%Li%NULL%1,         Yong%Gao%NULL%1,         Junjie%Zhou%NULL%1,         Sihua%Wang%NULL%1,         Bo%Huang%NULL%1,         Junhua%Wu%NULL%1,         Qin%Cao%NULL%1,         Yajun%Chen%NULL%1,         Zhenkao%Wang%NULL%1,         Danju%Luo%NULL%1,         Ting%Zhou%NULL%1,         Ruiting%Li%NULL%0,         You%Shang%NULL%0,         Xiu%Nie%whunhpath@163.com%1]'
$ws.Range("E10").Value = '[Wanyin%Tao%NULL%1,         Guorong%Zhang%NULL%1,         Xiaofang%Wang%NULL%1,         Meng%Guo%NULL%1,         Weihong%Zeng%NULL%1,         Zhihao%Xu%NULL%1,         Dan%Cao%NULL%1,         Aijun%Pan%NULL%0,         Yucai%Wang%NULL%1,         Kaiguang%Zhang%NULL%1,         Xiaoling%Ma%NULL%1,         Zhengxu%Chen%NULL%1,         Tengchuan%Jin%NULL%1,         Lianxin%Liu%NULL%0,         Jianping%Weng%NULL%1,         Shu%Zhu%NULL%1]'
$ws.Range("E11").Value = '[ L.%Yu%null%1,          Y.% Tong%null%1,          G.% Shen%null%1,          A.% Fu%null%1,          Y.% Lai%null%1,          X.% Zhou%null%2,          Y.% Yuan%null%1,          Y.% Wang%null%1,          Y.% Pan%null%1,          Z.% Yu%null%1,          Y.% Li%null%0,          T.% Liu%null%1,          H. % Jiang%null%1]'
$ws.Range("E12").Value = '[ J.-W.%Ai%null%1,          H.-C.% Zhang%null%1,          T.% Xu%null%1,          J.% Wu%null%1,          M.% Zhu%null%1,          Y.-Q.% Yu%null%1,          H.-Y.% Zhang%null%1,          Z.% Shen%null%1,          Y.% Li%null%0,          X.% Zhou%null%0,          G.-Q.% Zang%null%1,          J.% Xu%null%1,          W.-J.% Chen%null%1,          Y.-J.% Li%null%1,          D.-S.% Xie%null%1,          M.-Z.% Zhou%null%1,          J.-Y.% Sun%null%1,          J.-Z.% Chen%null%1,          W.-H. % Zhang%null%1]'
